$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Sheet1: update column C values (REST/window boundaries) and trim two rows ---
$newC = @(895, 1160, 2017, 2184, 2350, 2515, 3275)
for ($i = 0; $i -lt $newC.Length; $i++) {
    $ws1.Cells.Item(2 + $i, 3).Value = $newC[$i]
}
$ws1.Range("C9:C10").ClearContents()

# Move the active selection to match the updated REST window range
$ws1.Range("C4:C8").Select()

# --- Add Sheet2 (after Sheet1) with the electrode-data classification table ---
$ws2 = $wb.Worksheets.Add([Type]::Missing, $ws1)

$headers = @("FILENAMES", "FLEXION", "EXTENSION", "SUSTAIN", "REST", "WINDOW")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $c = $ws2.Cells.Item(1, 2 + $i)
    $c.Value = $headers[$i]
    $c.Font.Bold = $true
    $c.HorizontalAlignment = -4108
    $c.VerticalAlignment = -4160
    $c.Borders.LineStyle = 1
    $c.Borders.ColorIndex = -4105
}

$ws2.Range("B2").Value = "../data/electrodeData/CoolTerm Capture 2023-03-28 15-39-11.txt"

$dataRows = @(
    @(668, 767, 1792, 1990, 90),
    @(866, 965, 1891, 2056, $null),
    @(1064, 1164, $null, $null, $null),
    @(1263, 1395, $null, $null, $null),
    @(1527, 1659, $null, $null, $null),
    @(2188, 2287, $null, $null, $null),
    @(2386, 2486, $null, $null, $null)
)

for ($r = 0; $r -lt $dataRows.Length; $r++) {
    $row = 2 + $r
    $vals = $dataRows[$r]
    for ($cIdx = 0; $cIdx -lt $vals.Length; $cIdx++) {
        if ($null -ne $vals[$cIdx]) {
            $ws2.Cells.Item($row, 3 + $cIdx).Value = $vals[$cIdx]
        }
    }
}

# Column A: sequential index 0..17, rows 2..19, bold/centered/top/thin-bordered
for ($r = 0; $r -le 17; $r++) {
    $row = 2 + $r
    $c = $ws2.Cells.Item($row, 1)
    $c.Value = $r
    $c.Font.Bold = $true
    $c.HorizontalAlignment = -4108
    $c.VerticalAlignment = -4160
    $c.Borders.LineStyle = 1
    $c.Borders.ColorIndex = -4105
}

$ws1.Select()
